$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.888.27'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '2.431.43'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.41'
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.91'
$ws.Range("E6").Value = '  +2.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.562'
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '2.479.08'
$ws.Range("E9").Value = '  +1.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0988'
$ws.Range("E10").Value = '  +1.36%  '

$ws.Range("E11").Value = '  -1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.326'
$ws.Range("E12").Value = '  +1.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.70'
$ws.Range("E13").Value = '  -4.59%  '

$ws.Range("D14").Value = '2.874.55'
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").Value = '57.747.36'
$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.16'
$ws.Range("E16").Value = '  +2.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  +2.14%  '

$ws.Range("D18").Value = '2.440.21'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.44'
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.17'
$ws.Range("E20").Value = '  +1.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '317.51'
$ws.Range("E21").Value = '  +1.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("E22").Value = '  +6.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.72'
$ws.Range("E24").Value = '  -2.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.74'
$ws.Range("E25").Value = '  +1.45%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.992'
$ws.Range("E26").Value = '  -0.84%  '

$ws.Range("D27").Value = '2.535.55'
$ws.Range("E27").Value = '  -0.77%  '

$ws.Range("E28").Value = '  -3.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.155'
$ws.Range("E29").Value = '  -1.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.71'
$ws.Range("E30").Value = '  +6.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.24'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '0.0₃0745'
$ws.Range("E32").Value = '  +1.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.70'
$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.28'
$ws.Range("E34").Value = '  +2.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.16'
$ws.Range("E35").Value = '  +2.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.22%  '

$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.19'
$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.27'
$ws.Range("E39").Value = '  +7.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.91'
$ws.Range("E40").Value = '  +4.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.827'
$ws.Range("E41").Value = '  +4.69%  '

$ws.Range("E42").Value = '  +2.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.49'
$ws.Range("E43").Value = '  +0.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '136.27'
$ws.Range("E44").Value = '  +12.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.45'
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.08'
$ws.Range("E46").Value = '  +5.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '264.44'
$ws.Range("E47").Value = '  +0.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.581'
$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0503'
$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0217'
$ws.Range("E51").Value = '  +3.18%  '
